# Modification de la doc utilisateur:
#   - ajoute un paragraphe de texte "Test mathieu"
#   - ajoute un second paragraphe (vide) contenant le signet cache
#     "_GoBack" que Word depose automatiquement a la derniere position
#     d'edition lors de l'enregistrement.

$d = $word.ActiveDocument

# Le document de depart ne contient qu'un unique paragraphe vide ;
# on y ecrit le texte attendu.
$r = $d.Content
$r.InsertAfter("Test mathieu")

# On se place juste apres le texte insere pour y ajouter le nouveau
# paragraphe (vide) qui accueillera le bookmark "_GoBack".
$r.Collapse(0)

# InsertXML permet d'inserer le paragraphe + le bookmark directement,
# sans laisser de <w:r/> vide parasite dans le paragraphe final (ce que
# ferait InsertParagraphAfter + Bookmarks.Add).
$goBackXml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
  <pkg:part pkg:name="/_rels/.rels" pkg:contentType="application/vnd.openxmlformats-package.relationships+xml" pkg:padding="512">
    <pkg:xmlData>
      <Relationships xmlns="http://schemas.openxmlformats.org/package/2006/relationships">
        <Relationship Id="rId1" Type="http://schemas.openxmlformats.org/officeDocument/2006/relationships/officeDocument" Target="word/document.xml"/>
      </Relationships>
    </pkg:xmlData>
  </pkg:part>
  <pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
    <pkg:xmlData>
      <w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
        <w:body>
          <w:p>
            <w:bookmarkStart w:id="0" w:name="_GoBack"/>
            <w:bookmarkEnd w:id="0"/>
          </w:p>
        </w:body>
      </w:document>
    </pkg:xmlData>
  </pkg:part>
</pkg:package>
'@

[void]$r.InsertXML($goBackXml)
